$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: new "target" column label in O1
$ws.Range("O1").Value = "target"

# Data values (columns A-O, rows 2-18) replaced per diff
$ws.Range("A2").Value = 21
$ws.Range("B2").Value = 1.088576992162783
$ws.Range("C2").Value = 1.113879672048004
$ws.Range("D2").Value = 1.391852302160855
$ws.Range("E2").Value = 1.097675066956605
$ws.Range("F2").Value = 0.6193999948385083
$ws.Range("G2").Value = 1.260478162207732
$ws.Range("H2").Value = 0.2829759842181138
$ws.Range("I2").Value = -0.2048732487660765
$ws.Range("J2").Value = 1.280068775920671
$ws.Range("K2").Value = 0.4598623887690349
$ws.Range("L2").Value = 0.7524404323814086
$ws.Range("M2").Value = 0.802791065925304
$ws.Range("N2").Value = 0.7880756782855436
$ws.Range("O2").Value = 1
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = -1.052360446
$ws.Range("C3").Value = -1.116420485
$ws.Range("D3").Value = -1.102244233
$ws.Range("E3").Value = -2.994605003
$ws.Range("F3").Value = -2.034215286
$ws.Range("G3").Value = -1.796524334
$ws.Range("H3").Value = -1.041387772
$ws.Range("I3").Value = -0.212968827
$ws.Range("J3").Value = 0.294011831
$ws.Range("K3").Value = 0.103498634
$ws.Range("L3").Value = -0.09137593400000001
$ws.Range("M3").Value = -0.624227593
$ws.Range("N3").Value = -0.321369264
$ws.Range("O3").Value = 1
$ws.Range("A4").Value = 35
$ws.Range("B4").Value = -0.9673902202431353
$ws.Range("C4").Value = -1.005585475346329
$ws.Range("D4").Value = 0.01336928259588351
$ws.Range("E4").Value = -0.03144927010491952
$ws.Range("F4").Value = -0.8290864447953518
$ws.Range("G4").Value = -0.3559431232112886
$ws.Range("H4").Value = -1.188186603427425
$ws.Range("I4").Value = -2.503584584696504
$ws.Range("J4").Value = -2.388757987935037
$ws.Range("K4").Value = -1.560003476310523
$ws.Range("L4").Value = -0.9705462612823146
$ws.Range("M4").Value = -0.4379119648592557
$ws.Range("N4").Value = -0.6418788257929507
$ws.Range("O4").Value = 1
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = 1.09260717879037
$ws.Range("C5").Value = 1.10528016275872
$ws.Range("D5").Value = -0.5052371728717996
$ws.Range("E5").Value = -1.213035656723745
$ws.Range("F5").Value = -0.9397964201281946
$ws.Range("G5").Value = -0.07632542278754005
$ws.Range("H5").Value = -0.1814666348577954
$ws.Range("I5").Value = -0.08897598996320968
$ws.Range("J5").Value = -0.1461574093903557
$ws.Range("K5").Value = -0.5692170635129519
$ws.Range("L5").Value = -1.400241450620395
$ws.Range("M5").Value = -1.52723548122241
$ws.Range("N5").Value = -1.861333063595394
$ws.Range("O5").Value = 1
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = 1.168895327646786
$ws.Range("C6").Value = 1.124790058217134
$ws.Range("D6").Value = -0.108837563652328
$ws.Range("E6").Value = -0.527461071075682
$ws.Range("F6").Value = -0.8652942322810397
$ws.Range("G6").Value = 0.1168204941201648
$ws.Range("H6").Value = -0.723641543470295
$ws.Range("I6").Value = -0.498731630075758
$ws.Range("J6").Value = -1.025476842083122
$ws.Range("K6").Value = -1.622354741293647
$ws.Range("L6").Value = -2.117899344682976
$ws.Range("M6").Value = -2.105917768685274
$ws.Range("N6").Value = -2.604525295334807
$ws.Range("O6").Value = 1
$ws.Range("A7").Value = 26
$ws.Range("B7").Value = -0.9801136137722205
$ws.Range("C7").Value = -1.041956232283092
$ws.Range("D7").Value = -0.1886708569035639
$ws.Range("E7").Value = 0.03187326138358325
$ws.Range("F7").Value = 0.02696679057833097
$ws.Range("G7").Value = 0.4139021714331789
$ws.Range("H7").Value = -0.2537263285459994
$ws.Range("I7").Value = -0.7238580530731326
$ws.Range("J7").Value = -0.8867859043564654
$ws.Range("K7").Value = -0.4421644275965932
$ws.Range("L7").Value = 0.5141295735833261
$ws.Range("M7").Value = 0.4309079702258571
$ws.Range("N7").Value = 0.6473643899869618
$ws.Range("O7").Value = 1
$ws.Range("A8").Value = 8
$ws.Range("B8").Value = 1.200056980136892
$ws.Range("C8").Value = 1.123243570411099
$ws.Range("D8").Value = -0.55241059474163
$ws.Range("E8").Value = -1.036721346042088
$ws.Range("F8").Value = -0.9308365164426041
$ws.Range("G8").Value = 0.1155780199142884
$ws.Range("H8").Value = 0.2663936930680734
$ws.Range("I8").Value = 0.02870485741260834
$ws.Range("J8").Value = -0.602788081833072
$ws.Range("K8").Value = -1.192180913668596
$ws.Range("L8").Value = -2.247474818272565
$ws.Range("M8").Value = -2.537178146746713
$ws.Range("N8").Value = -2.569831016425473
$ws.Range("O8").Value = 1
$ws.Range("A9").Value = 5
$ws.Range("B9").Value = 0.09099733199999999
$ws.Range("C9").Value = 0.165123739
$ws.Range("D9").Value = 1.033777922
$ws.Range("E9").Value = 0.62781765
$ws.Range("F9").Value = 0.7088679729999999
$ws.Range("G9").Value = 0.9639872290000001
$ws.Range("H9").Value = -0.570432089
$ws.Range("I9").Value = -0.6886823870000001
$ws.Range("J9").Value = -0.226930379
$ws.Range("K9").Value = 0.2216372
$ws.Range("L9").Value = -0.154145444
$ws.Range("M9").Value = 0.244069248
$ws.Range("N9").Value = 0.37543717
$ws.Range("O9").Value = 1
$ws.Range("A10").Value = 19
$ws.Range("B10").Value = -1.053461550083322
$ws.Range("C10").Value = -1.013951085727377
$ws.Range("D10").Value = -1.670975857646011
$ws.Range("E10").Value = -2.305474877223224
$ws.Range("F10").Value = -1.061541930846145
$ws.Range("G10").Value = -1.097802694772028
$ws.Range("H10").Value = -0.2332301376733607
$ws.Range("I10").Value = 0.2441601328574322
$ws.Range("J10").Value = 0.3523387254212044
$ws.Range("K10").Value = -0.2489324891736348
$ws.Range("L10").Value = -0.502001930595651
$ws.Range("M10").Value = -0.6653240106200277
$ws.Range("N10").Value = -0.3590114128251046
$ws.Range("O10").Value = 1
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 0.08650412199999999
$ws.Range("C11").Value = 0.16076415
$ws.Range("D11").Value = 1.287529653
$ws.Range("E11").Value = 1.458507674
$ws.Range("F11").Value = 1.708264452
$ws.Range("G11").Value = 1.737302749
$ws.Range("H11").Value = 1.454654867
$ws.Range("I11").Value = 1.685926169
$ws.Range("J11").Value = 1.198936477
$ws.Range("K11").Value = 1.009485194
$ws.Range("L11").Value = 0.6228660539999999
$ws.Range("M11").Value = 0.391809069
$ws.Range("N11").Value = 0.180105296
$ws.Range("O11").Value = 2
$ws.Range("A12").Value = 47
$ws.Range("B12").Value = -0.915313197917203
$ws.Range("C12").Value = -0.8002599253251903
$ws.Range("D12").Value = -0.4858008915469419
$ws.Range("E12").Value = -0.459102613136871
$ws.Range("F12").Value = -0.3886572214119309
$ws.Range("G12").Value = -1.131570294461621
$ws.Range("H12").Value = -0.4310121901661068
$ws.Range("I12").Value = 0.1352166989295445
$ws.Range("J12").Value = 0.0321162105799529
$ws.Range("K12").Value = 0.2689214043482228
$ws.Range("L12").Value = 0.1719394788427852
$ws.Range("M12").Value = 0.01614742187502364
$ws.Range("N12").Value = 0.2079486552789812
$ws.Range("O12").Value = 2
$ws.Range("A13").Value = 39
$ws.Range("B13").Value = -0.5773369837013437
$ws.Range("C13").Value = -0.7008086824709544
$ws.Range("D13").Value = -1.04388155837972
$ws.Range("E13").Value = -0.889075565958088
$ws.Range("F13").Value = -0.4777298958116295
$ws.Range("G13").Value = -0.7205425914449428
$ws.Range("H13").Value = -0.7971610748256895
$ws.Range("I13").Value = -0.008196996251951846
$ws.Range("J13").Value = -1.049105483397343
$ws.Range("K13").Value = -0.1228300179444196
$ws.Range("L13").Value = 0.0830021749573447
$ws.Range("M13").Value = 0.2045357481978724
$ws.Range("N13").Value = 0.269263486568356
$ws.Range("O13").Value = 2
$ws.Range("A14").Value = 51
$ws.Range("B14").Value = -0.8999088371623735
$ws.Range("C14").Value = -0.7979083676895516
$ws.Range("D14").Value = -0.3777293800015497
$ws.Range("E14").Value = -0.4738492208868103
$ws.Range("F14").Value = -0.4220898121785666
$ws.Range("G14").Value = -0.6004718027471426
$ws.Range("H14").Value = -1.011272578681248
$ws.Range("I14").Value = 0.1028459493065351
$ws.Range("J14").Value = 0.3114022967378419
$ws.Range("K14").Value = 0.1371577873176834
$ws.Range("L14").Value = -0.03766079638830301
$ws.Range("M14").Value = 0.1002645347072279
$ws.Range("N14").Value = 0.2894954110783645
$ws.Range("O14").Value = 2
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = -0.9819434020000001
$ws.Range("C15").Value = -1.027168868
$ws.Range("D15").Value = -0.391420759
$ws.Range("E15").Value = -1.071957494
$ws.Range("F15").Value = -1.015058936
$ws.Range("G15").Value = -1.251488744
$ws.Range("H15").Value = -1.889056123
$ws.Range("I15").Value = -0.649593132
$ws.Range("J15").Value = -0.321722455
$ws.Range("K15").Value = -0.310259495
$ws.Range("L15").Value = -0.362302266
$ws.Range("M15").Value = -0.24811692
$ws.Range("N15").Value = 0.142904433
$ws.Range("O15").Value = 2
$ws.Range("A16").Value = 44
$ws.Range("B16").Value = -0.4486676309063777
$ws.Range("C16").Value = -0.6981368151440918
$ws.Range("D16").Value = -1.893864632384399
$ws.Range("E16").Value = -0.8138407903746437
$ws.Range("F16").Value = -1.131585928047311
$ws.Range("G16").Value = -1.521676146207622
$ws.Range("H16").Value = -0.421539501489349
$ws.Range("I16").Value = -0.9139190013134963
$ws.Range("J16").Value = 0.078154540762746
$ws.Range("K16").Value = -0.03860679876158272
$ws.Range("L16").Value = -0.2122522625214683
$ws.Range("M16").Value = -0.09766339219853593
$ws.Range("N16").Value = 0.07815511282505004
$ws.Range("O16").Value = 2
$ws.Range("A17").Value = 59
$ws.Range("B17").Value = 1.098824281425039
$ws.Range("C17").Value = 1.124544984094857
$ws.Range("D17").Value = 1.730050599012491
$ws.Range("E17").Value = 1.85323807668174
$ws.Range("F17").Value = 2.19314676580624
$ws.Range("G17").Value = 2.024663755796917
$ws.Range("H17").Value = 2.554541131790494
$ws.Range("I17").Value = 2.318109608543673
$ws.Range("J17").Value = 2.100052046017395
$ws.Range("K17").Value = 2.044557233732532
$ws.Range("L17").Value = 1.594573606845447
$ws.Range("M17").Value = 1.329427998478601
$ws.Range("N17").Value = 1.092029141150236
$ws.Range("O17").Value = 3
$ws.Range("A18").Value = 63
$ws.Range("B18").Value = -0.6057004064306782
$ws.Range("C18").Value = -0.5377353754033299
$ws.Range("D18").Value = 0.06023393650398268
$ws.Range("E18").Value = 1.156450741882419
$ws.Range("F18").Value = 1.695659991481
$ws.Range("G18").Value = 0.5459799358931066
$ws.Range("H18").Value = 1.737616217432471
$ws.Range("I18").Value = 1.459572815406067
$ws.Range("J18").Value = 1.27574394030813
$ws.Range("K18").Value = 1.463114216622488
$ws.Range("L18").Value = 1.41933957798711
$ws.Range("M18").Value = 1.227404464743771
$ws.Range("N18").Value = 0.9388314252113402
$ws.Range("O18").Value = 3
